$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (Excel would otherwise auto-coerce strings
# like "1.003" or "314.67" into numbers). Force the cell to text format,
# assign, then reset the style back to Normal so no stray formatting is
# left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "27.638.37"
$ws.Range("E2").Value = "  -2.37%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.846.11"
$ws.Range("E3").Value = "  -1.03%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "314.67"
$ws.Range("E5").Value = "  -1.60%  "

# Row 6 - USDC
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.09%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -2.95%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3659"
$ws.Range("E8").Value = "  -2.22%  "

# Row 9 - OKB
Set-TextValue $ws.Range("D9") "46.01"
$ws.Range("E9").Value = "  +1.85%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.07253"
$ws.Range("E10").Value = "  -4.17%  "

# Row 11 - Polygon
Set-TextValue $ws.Range("D11") "0.9002"
$ws.Range("E11").Value = "  -4.33%  "

# Row 12 - Solana
Set-TextValue $ws.Range("D12") "20.69"
$ws.Range("E12").Value = "  -2.98%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.822.06"
$ws.Range("E13").Value = "  -2.91%  "

# Row 14/15 - Polkadot and Chainlink swap places
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "6.573"
$ws.Range("E14").Value = "  -2.43%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "5.360"
$ws.Range("E15").Value = "  -2.35%  "

# Row 16 - TRON
Set-TextValue $ws.Range("D16") "0.06833"
$ws.Range("E16").Value = "  -0.42%  "

# Row 17 - BinanceUSD
Set-TextValue $ws.Range("D17") "1.004"
$ws.Range("E17").Value = "  -0.06%  "

# Row 18 - Litecoin
Set-TextValue $ws.Range("D18") "77.90"
$ws.Range("E18").Value = "  -5.25%  "

# Row 19 - ShibaInu
Set-TextValue $ws.Range("D19") "0.000008817"
$ws.Range("E19").Value = "  -3.28%  "

# Row 20 - Dai
Set-TextValue $ws.Range("D20") "1.001"
$ws.Range("E20").Value = "  -0.11%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  -3.79%  "

# Row 22 - WrappedBTC
Set-TextValue $ws.Range("D22") "27.624.33"
$ws.Range("E22").Value = "  -2.41%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "4.961"
$ws.Range("E23").Value = "  -3.91%  "

# Row 24 - Cosmos
Set-TextValue $ws.Range("D24") "10.62"
$ws.Range("E24").Value = "  -1.41%  "

# Row 25 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D25") "2.085.85"
$ws.Range("E25").Value = "  +0.01%  "

# Row 26 - Toncoin
Set-TextValue $ws.Range("D26") "2.047"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27 - Monero
Set-TextValue $ws.Range("D27") "154.00"
$ws.Range("E27").Value = "  -0.54%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "18.23"
$ws.Range("E28").Value = "  -1.05%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "5.291"
$ws.Range("E29").Value = "  -1.64%  "

# Row 30 - LidoDAOToken
Set-TextValue $ws.Range("D30") "1.819"
$ws.Range("E30").Value = "  +4.78%  "

# Row 31 - BitcoinCash
Set-TextValue $ws.Range("D31") "110.79"
$ws.Range("E31").Value = "  -3.50%  "

# Row 32 - Stellar
Set-TextValue $ws.Range("D32") "0.08884"
$ws.Range("E32").Value = "  -1.99%  "

# Row 33 - ImmutableX
Set-TextValue $ws.Range("D33") "0.7701"
$ws.Range("E33").Value = "  -4.68%  "

# Row 34 - Filecoin
Set-TextValue $ws.Range("D34") "4.544"

# Row 35 - HuobiToken
Set-TextValue $ws.Range("D35") "2.974"
$ws.Range("E35").Value = "  +0.90%  "

# Row 36 - ARBITRUM
Set-TextValue $ws.Range("D36") "1.084"
$ws.Range("E36").Value = "  -8.03%  "

# Row 37 - Frax
Set-TextValue $ws.Range("D37") "1.000"
$ws.Range("E37").Value = "  -0.17%  "

# Row 38/39 - Hedera and TrustWalletToken swap places
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.099"
$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.05390"
$ws.Range("E39").Value = "  -1.82%  "

# Row 40/41 - VeChain and MXToken swap places
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.956"
$ws.Range("E40").Value = "  -1.70%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.01925"
$ws.Range("E41").Value = "  -1.78%  "

# Row 42 - TheSandbox
Set-TextValue $ws.Range("D42") "0.5065"
$ws.Range("E42").Value = "  -3.97%  "

# Row 43 - FraxShare
Set-TextValue $ws.Range("D43") "6.806"
$ws.Range("E43").Value = "  -5.33%  "

# Row 44 - Algorand
Set-TextValue $ws.Range("D44") "0.1643"
$ws.Range("E44").Value = "  -2.25%  "

# Row 45 - Aptos
Set-TextValue $ws.Range("D45") "8.232"
$ws.Range("E45").Value = "  -7.05%  "

# Row 46 - Cronos
Set-TextValue $ws.Range("D46") "0.06641"
$ws.Range("E46").Value = "  -2.15%  "

# Row 47 - EnergySwap
Set-TextValue $ws.Range("D47") "10.37"
$ws.Range("E47").Value = "  -2.01%  "

# Row 48 - Decentraland
Set-TextValue $ws.Range("D48") "0.4725"
$ws.Range("E48").Value = "  -3.67%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "105.26"

# Row 50 - PaxDollar
Set-TextValue $ws.Range("D50") "1.001"
$ws.Range("E50").Value = "  -0.11%  "

# Row 51 - NEARProtocol
Set-TextValue $ws.Range("D51") "1.638"
$ws.Range("E51").Value = "  -3.27%  "
